$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap full row contents (country name + stats) for the three reordered
# country pairs in the "paises.xlsx" source (Santa Lucia/Belice,
# Namibia/San Vicente y las Granadinas, Burundi/San Cristobal y Nieves). ---

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("A$rowA`:H$rowA")
    $rangeB = $ws.Range("A$rowB`:H$rowB")
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

Swap-Rows 188 189
Swap-Rows 194 195
Swap-Rows 198 199

# --- Updated daily statistics for a handful of countries ---

# Pakistan (row 27)
$ws.Range("B27").Value = 20186
$ws.Range("C27").Value = 102
$ws.Range("D27").Value = 5590
$ws.Range("E27").Value = 14134
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 462

# Kazajistan (row 60)
$ws.Range("B60").Value = 3964
$ws.Range("C60").Value = 44
$ws.Range("E60").Value = 2853

# Hungria (row 63)
$ws.Range("B63").Value = 3035
$ws.Range("C63").Value = 37
$ws.Range("D63").Value = 630
$ws.Range("E63").Value = 2054
$ws.Range("F63").Value = 55
$ws.Range("G63").Value = 11
$ws.Range("H63").Value = 351

# Tailandia (row 64)
$ws.Range("B64").Value = 2987
$ws.Range("C64").Value = 18
$ws.Range("D64").Value = 2740
$ws.Range("E64").Value = 193

# Kirguistan (row 96)
$ws.Range("F96").Value = 13

Write-Host "Done applying country/provincia updates"
